$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1268
$ws.Range("F5").Value = 31
$ws.Range("F7").Value = 1423
$ws.Range("F8").Value = 81
$ws.Range("F9").Value = 30
$ws.Range("F10").Value = 674
$ws.Range("F11").Value = 149
$ws.Range("F12").Value = 158
$ws.Range("F14").Value = 1311
$ws.Range("F15").Value = 11
$ws.Range("F16").Value = 517
$ws.Range("F20").Value = 119
$ws.Range("F21").Value = 748
$ws.Range("F22").Value = 2556
$ws.Range("F28").Value = 11
$ws.Range("F29").Value = 121
$ws.Range("F31").Value = 929
$ws.Range("F33").Value = 116

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 640
$ws.Range("F9").Value = 283
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 507
$ws.Range("F16").Value = 944
$ws.Range("G26").Value = 90

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2318
$ws.Range("F6").Value = 926
$ws.Range("F9").Value = 1170
$ws.Range("F10").Value = 291

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2318
$ws.Range("F8").Value = 926
$ws.Range("F9").Value = 1170
$ws.Range("F10").Value = 291
$ws.Range("F12").Value = 1268
$ws.Range("F13").Value = 31
$ws.Range("F15").Value = 1423
$ws.Range("F16").Value = 81
$ws.Range("F17").Value = 30
$ws.Range("F18").Value = 674
$ws.Range("F19").Value = 149
$ws.Range("F21").Value = 158
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 517
$ws.Range("F26").Value = 119
$ws.Range("F27").Value = 748
$ws.Range("F28").Value = 2556
$ws.Range("F32").Value = 283
$ws.Range("F34").Value = 121
$ws.Range("F36").Value = 929
$ws.Range("F37").Value = 507
$ws.Range("F40").Value = 116
